$d = $word.ActiveDocument

# 1) Update the intro sentence above the (soon to be) table.
$d.Content.Find.Execute(
    "To support decision-making when resource conflicts arise, project priorities are defined in 3 dimensions:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The table below defines the project priorities across three dimensions to facilitate decision-making:",
    2
) | Out-Null

$introPara = 144

# 2) Delete the 8 bullet paragraphs (Constraints/Schedule/.../Scope) that follow.
$firstBullet = $d.Paragraphs($introPara + 1)
$lastBullet = $d.Paragraphs($introPara + 8)
$deleteRange = $d.Range($firstBullet.Range.Start, $lastBullet.Range.End)
$deleteRange.Delete() | Out-Null

# 3) Insert the replacement table at the now-collapsed insertion point.
$insertAt = $d.Paragraphs($introPara + 1).Range.Start
$insertRange = $d.Range($insertAt, $insertAt)
$tableXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/><w:gridCol w:w="2160"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Dimension</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Constraint (Must adhere to)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Driver (Key success factor)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Degree of Freedom (Adjustable)</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Schedule</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Must complete Release 1.0 before [End Date] for acceptance.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Budget</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Implementation cost ~0 (Using student manpower, free tier).</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Quality</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Data accuracy is paramount. System stability.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Usability</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Extremely friendly interface for non-IT staff.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Scope</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="2160"/></w:tcPr><w:p><w:r><w:t>Advanced features (Mobile App, Accounting API) can be cut or delayed to Release 2.0.</w:t></w:r></w:p></w:tc></w:tr></w:tbl></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($tableXml) | Out-Null

# 4) Add a single blank paragraph between the table and the following
#    Heading2 ("3.3. Deployment Considerations"). InsertXML at this exact
#    table-end boundary needs two <w:p/> nodes to materialize one new
#    paragraph, so insert two and then delete the extra one.
$newTbl = $d.Tables($d.Tables.Count)
$tblEnd = $newTbl.Range.End
$afterTableRange = $d.Range($tblEnd, $tblEnd)
$paraXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$afterTableRange.InsertXML($paraXml) | Out-Null
$extraBlank = $d.Range($tblEnd, $tblEnd + 1)
$extraBlank.Delete() | Out-Null

Write-Output "done"
